$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply the AutoFilter on column A (Modulo) so only the selected values remain
# visible: "Viaje -> Operador", "Viaje->Cobros" and "Viajes". This hides every
# other data row (2-20) and keeps rows 21-25 visible.
[void]$ws.Range("A1:D25").AutoFilter(1, @("Viaje -> Operador","Viaje->Cobros","Viajes"), 7)

# Highlight the relevant rows in yellow (adds a new fill/style), skipping the
# header row and row 23.
$ws.Range("A2:D22").Interior.Color = 65535
$ws.Range("A24:D25").Interior.Color = 65535

# Move the active selection to D23.
[void]$ws.Range("D23").Select()
